# Censuses DB adapted to MySQL
#
# The voting-station flag in column D used to be stored as text ("A"/"B",
# shared-string cells). Moving the census backend to MySQL, that column is
# now a plain numeric flag, so every row gets a literal 1 instead of the
# old letter codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D2/D3: "A"/"B" (text) -> 1 (number)
$ws.Range("D2:D3").Value = 1

# The active selection in the saved sheet moved from F2 to D9.
$ws.Activate()
$ws.Range("D9").Select()

# Cosmetic/localization touch-ups from the same commit (Spanish Excel UI):
# the built-in "Hyperlink" cell style is shown as "Hipervínculo". Best
# effort only - older/limited hosts may not expose a writable style name,
# so this is guarded to avoid breaking the rest of the script.
try {
    $wb.Styles.Item("Hyperlink").Name = "Hipervínculo"
} catch {
}
